$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63-164 down to 64-165.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the new record's data.
$ws.Cells.Item(63, 1).Value2 = 10
$ws.Cells.Item(63, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value2 = "La Araucanía"
$ws.Cells.Item(63, 4).Value2 = 44533
$ws.Cells.Item(63, 5).Value2 = 9
$ws.Cells.Item(63, 6).Value2 = 100112005
$ws.Cells.Item(63, 7).Value2 = "Puerro"
$ws.Cells.Item(63, 8).Value2 = "Azul de Maquehue"
$ws.Cells.Item(63, 9).Value2 = "Primera"
$ws.Cells.Item(63, 10).Value2 = 45
$ws.Cells.Item(63, 11).Value2 = 8000
$ws.Cells.Item(63, 12).Value2 = 8000
$ws.Cells.Item(63, 13).Value2 = 8000
$ws.Cells.Item(63, 14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(63, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(63, 16).Value2 = 667
$ws.Cells.Item(63, 17).Value2 = 12
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"
